# Adjust Investment Summary table column widths for better formatting
# (also tightens the Timeline & Milestones table's columns to match)
#
# PowerPoint COM reports/accepts table Column.Width in points; OOXML
# stores the grid in EMU (1 pt = 12700 EMU). The values below are the
# exact EMU targets converted to points so the round-trip is lossless.

$p = $ppt.ActivePresentation

function Find-SlideByTitle($pres, $title) {
    for ($i = 1; $i -le $pres.Slides.Count; $i++) {
        $slide = $pres.Slides.Item($i)
        if ($slide.Shapes.Count -ge 1) {
            $txt = ""
            try { $txt = $slide.Shapes.Item(1).TextFrame.TextRange.Text } catch {}
            if ($txt -eq $title) {
                return $slide
            }
        }
    }
    return $null
}

function Find-TableShape($slide) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            return $shp
        }
    }
    return $null
}

# --- Timeline & Milestones table: 4 columns -------------------------------
$timelineSlide = Find-SlideByTitle $p "Timeline & Milestones"
$timelineShape = Find-TableShape $timelineSlide
$timelineTbl = $timelineShape.Table

$timelineTbl.Columns.Item(1).Width = (871093 / 12700)
$timelineTbl.Columns.Item(2).Width = (2177733 / 12700)
$timelineTbl.Columns.Item(3).Width = (1306639 / 12700)
$timelineTbl.Columns.Item(4).Width = (4355466 / 12700)

# --- Investment Summary table: 7 columns ----------------------------------
$investSlide = Find-SlideByTitle $p "Investment Summary"
$investShape = Find-TableShape $investSlide
$investTbl = $investShape.Table

$investTbl.Columns.Item(1).Width = (1742186 / 12700)
$investTbl.Columns.Item(2).Width = (1045311 / 12700)
$investTbl.Columns.Item(3).Width = (2003514 / 12700)
$investTbl.Columns.Item(4).Width = (1132421 / 12700)
$investTbl.Columns.Item(5).Width = (871093 / 12700)
$investTbl.Columns.Item(6).Width = (871093 / 12700)
$investTbl.Columns.Item(7).Width = (1045311 / 12700)
